$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 223
$ws1.Range("F3").Value = 1401
$ws1.Range("F4").Value = 19684
$ws1.Range("F5").Value = 790
$ws1.Range("F6").Value = 305
$ws1.Range("F7").Value = 1092
$ws1.Range("F8").Value = 0
$ws1.Range("F9").Value = 7453
$ws1.Range("F10").Value = 495
$ws1.Range("F11").Value = 728
$ws1.Range("F12").Value = 256
$ws1.Range("F13").Value = 35
$ws1.Range("F14").Value = 149
$ws1.Range("F15").Value = 106
$ws1.Range("F17").Value = 231
$ws1.Range("F18").Value = 187
$ws1.Range("F19").Value = 1331
$ws1.Range("F20").Value = 385
$ws1.Range("F21").Value = 70
$ws1.Range("F22").Value = 676
$ws1.Range("F23").Value = 46
$ws1.Range("F24").Value = 51
$ws1.Range("F25").Value = 61
$ws1.Range("F28").Value = 25
$ws1.Range("F29").Value = 12
$ws1.Range("F30").Value = 169
$ws1.Range("F31").Value = 5220
$ws1.Range("F33").Value = 50
$ws1.Range("F34").Value = 0
$ws1.Range("F36").Value = 86
$ws1.Range("F38").Value = 12540
$ws1.Range("F39").Value = 1324
$ws1.Range("F40").Value = 63
$ws1.Range("F41").Value = 18
$ws1.Range("F42").Value = 54
$ws1.Range("F43").Value = 252
$ws1.Range("F44").Value = 346
$ws1.Range("F45").Value = 3979
$ws1.Range("F46").Value = 0
$ws1.Range("F47").Value = 93

# Sheet 2: 演出 (Performances)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F2").Value = 2

# Sheet 4: 全部类型 (All types)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 223
$ws4.Range("F3").Value = 1401
$ws4.Range("F4").Value = 0
$ws4.Range("F5").Value = 790
$ws4.Range("F6").Value = 305
$ws4.Range("F7").Value = 1092
$ws4.Range("F9").Value = 7453
$ws4.Range("F10").Value = 0
$ws4.Range("F12").Value = 256
$ws4.Range("F13").Value = 35
$ws4.Range("F14").Value = 149
$ws4.Range("F15").Value = 106
$ws4.Range("F16").Value = 2
$ws4.Range("F17").Value = 0
$ws4.Range("F18").Value = 187
$ws4.Range("F19").Value = 1331
$ws4.Range("F20").Value = 385
$ws4.Range("F21").Value = 70
$ws4.Range("F25").Value = 61
$ws4.Range("F26").Value = 314
$ws4.Range("F27").Value = 1080
$ws4.Range("F28").Value = 25
$ws4.Range("F29").Value = 12
$ws4.Range("F30").Value = 169
$ws4.Range("F32").Value = 556
$ws4.Range("F34").Value = 0
$ws4.Range("F35").Value = 33
$ws4.Range("F36").Value = 2784
$ws4.Range("F39").Value = 16
$ws4.Range("F43").Value = 18
$ws4.Range("F44").Value = 54
$ws4.Range("F45").Value = 252
$ws4.Range("F46").Value = 346
$ws4.Range("F47").Value = 3979
$ws4.Range("F49").Value = 0
